{"js": "// Applies the within-100 answer-table update: each old arithmetic\n// expression in the table is replaced with its new value (commit\n// \"Update master to output generated at aa3dc9e\").\nconst replacements = [\n  [\"5+48=53\", \"42-10=32\"],\n  [\"44-40=4\", \"63+13=76\"],\n  [\"15-12=3\", \"15+59=74\"],\n  [\"93+2=95\", \"90-6=84\"],\n  [\"77-51=26\", \"83+6=89\"],\n  [\"82-20=62\", \"52+30=82\"],\n  [\"14+1=15\", \"84-82=2\"],\n  [\"28+24=52\", \"7+45=52\"],\n  [\"3+18=21\", \"84-43=41\"],\n  [\"95-73=22\", \"71+3=74\"],\n  [\"30+37=67\", \"11+48=59\"],\n  [\"91-82=9\", \"27-19=8\"],\n  [\"36+3=39\", \"4+26=30\"],\n  [\"21+36=57\", \"56-29=27\"],\n  [\"8+8=16\", \"77-47=30\"],\n  [\"96-53=43\", \"89-11=78\"],\n  [\"5+17=22\", \"85-71=14\"],\n  [\"7-2=5\", \"89-54=35\"],\n  [\"32+22=54\", \"24+22=46\"],\n  [\"36+14=50\", \"94-37=57\"],\n  [\"73-46=27\", \"0+78=78\"],\n  [\"7+92=99\", \"33-6=27\"],\n  [\"6+30=36\", \"24+48=72\"],\n  [\"92-13=79\", \"1+57=58\"],\n  [\"91-12=79\", \"26+7=33\"],\n  [\"36+62=98\", \"84+7=91\"],\n  [\"13+2=15\", \"65+14=79\"],\n  [\"28+51=79\", \"17+69=86\"],\n  [\"57-9=48\", \"78-16=62\"],\n  [\"36+35=71\", \"67-6=61\"],\n  [\"59-19=40\", \"39-1=38\"],\n  [\"84-36=48\", \"10+79=89\"],\n  [\"37-20=17\", \"47+11=58\"],\n  [\"67+19=86\", \"94-75=19\"],\n  [\"30+54=84\", \"19+63=82\"],\n  [\"26+19=45\", \"10+32=42\"],\n  [\"0+14=14\", \"7+38=45\"],\n  [\"44-3=41\", \"91-76=15\"],\n  [\"50-13=37\", \"95-10=85\"],\n  [\"37+18=55\", \"32+61=93\"],\n  [\"4+53=57\", \"23-12=11\"],\n  [\"6+83=89\", \"57+21=78\"],\n  [\"4+6=10\", \"13+83=96\"],\n  [\"3+7=10\", \"95-20=75\"],\n  [\"37+47=84\", \"36-32=4\"],\n  [\"67-21=46\", \"63-26=37\"],\n  [\"97-37=60\", \"98-47=51\"],\n  [\"18+20=38\", \"92-48=44\"],\n  [\"89-6=83\", \"98-31=67\"],\n  [\"29+3=32\", \"10+85=95\"],\n  [\"38-11=27\", \"87-38=49\"],\n  [\"2+43=45\", \"26+21=47\"],\n  [\"16+72=88\", \"70-30=40\"],\n  [\"2+91=93\", \"23+0=23\"],\n  [\"52-28=24\", \"26-12=14\"],\n  [\"27+69=96\", \"90-56=34\"],\n  [\"45-38=7\", \"40-4=36\"],\n  [\"6+60=66\", \"38-3=35\"],\n  [\"58-18=40\", \"50-36=14\"],\n  [\"90-65=25\", \"25+40=65\"],\n  [\"64+28=92\", \"51+30=81\"],\n  [\"60+27=87\", \"76-25=51\"],\n  [\"80-53=27\", \"40+29=69\"],\n  [\"69-3=66\", \"10+89=99\"],\n  [\"1+39=40\", \"79-24=55\"],\n  [\"44+53=97\", \"85+13=98\"],\n  [\"36-34=2\", \"98-61=37\"],\n  [\"77+17=94\", \"55+29=84\"],\n  [\"11+22=33\", \"1+53=54\"],\n  [\"4+68=72\", \"1+17=18\"],\n  [\"51-34=17\", \"1+30=31\"],\n  [\"55-5=50\", \"85-84=1\"],\n  [\"63+30=93\", \"74+14=88\"],\n  [\"24+75=99\", \"42+51=93\"],\n  [\"44+52=96\", \"11+73=84\"],\n  [\"51-36=15\", \"49+44=93\"],\n  [\"11+76=87\", \"71-53=18\"],\n  [\"38+32=70\", \"26+8=34\"],\n  [\"68-47=21\", \"76-32=44\"],\n  [\"89-57=32\", \"44+15=59\"],\n  [\"48+36=84\", \"62+31=93\"],\n  [\"64-34=30\", \"43+35=78\"],\n  [\"73-67=6\", \"70-7=63\"],\n  [\"4+75=79\", \"76+1=77\"],\n  [\"57+42=99\", \"13+33=46\"],\n  [\"19+55=74\", \"7+85=92\"],\n  [\"28+18=46\", \"4+11=15\"],\n  [\"68-53=15\", \"35+43=78\"],\n  [\"96-27=69\", \"17-2=15\"],\n  [\"51-26=25\", \"2+45=47\"],\n  [\"19+54=73\", \"37-17=20\"],\n  [\"77-0=77\", \"64-10=54\"],\n  [\"80-72=8\", \"41+45=86\"],\n  [\"19-17=2\", \"54-2=52\"],\n  [\"23+76=99\", \"68-57=11\"],\n  [\"3+14=17\", \"94-70=24\"],\n  [\"41+22=63\", \"33+48=81\"],\n  [\"29+46=75\", \"62-19=43\"],\n  [\"43-28=15\", \"71+27=98\"],\n  [\"53+7=60\", \"98-96=2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Replace only the first occurrence (old values are unique in this\n  // document, so there should be exactly one match each).\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Applies the within-100 answer-table update: each old arithmetic\n# expression in the table is replaced with its new value (commit\n# \"Update master to output generated at aa3dc9e\").\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '5+48=53'; New = '42-10=32' },\n    @{ Old = '44-40=4'; New = '63+13=76' },\n    @{ Old = '15-12=3'; New = '15+59=74' },\n    @{ Old = '93+2=95'; New = '90-6=84' },\n    @{ Old = '77-51=26'; New = '83+6=89' },\n    @{ Old = '82-20=62'; New = '52+30=82' },\n    @{ Old = '14+1=15'; New = '84-82=2' },\n    @{ Old = '28+24=52'; New = '7+45=52' },\n    @{ Old = '3+18=21'; New = '84-43=41' },\n    @{ Old = '95-73=22'; New = '71+3=74' },\n    @{ Old = '30+37=67'; New = '11+48=59' },\n    @{ Old = '91-82=9'; New = '27-19=8' },\n    @{ Old = '36+3=39'; New = '4+26=30' },\n    @{ Old = '21+36=57'; New = '56-29=27' },\n    @{ Old = '8+8=16'; New = '77-47=30' },\n    @{ Old = '96-53=43'; New = '89-11=78' },\n    @{ Old = '5+17=22'; New = '85-71=14' },\n    @{ Old = '7-2=5'; New = '89-54=35' },\n    @{ Old = '32+22=54'; New = '24+22=46' },\n    @{ Old = '36+14=50'; New = '94-37=57' },\n    @{ Old = '73-46=27'; New = '0+78=78' },\n    @{ Old = '7+92=99'; New = '33-6=27' },\n    @{ Old = '6+30=36'; New = '24+48=72' },\n    @{ Old = '92-13=79'; New = '1+57=58' },\n    @{ Old = '91-12=79'; New = '26+7=33' },\n    @{ Old = '36+62=98'; New = '84+7=91' },\n    @{ Old = '13+2=15'; New = '65+14=79' },\n    @{ Old = '28+51=79'; New = '17+69=86' },\n    @{ Old = '57-9=48'; New = '78-16=62' },\n    @{ Old = '36+35=71'; New = '67-6=61' },\n    @{ Old = '59-19=40'; New = '39-1=38' },\n    @{ Old = '84-36=48'; New = '10+79=89' },\n    @{ Old = '37-20=17'; New = '47+11=58' },\n    @{ Old = '67+19=86'; New = '94-75=19' },\n    @{ Old = '30+54=84'; New = '19+63=82' },\n    @{ Old = '26+19=45'; New = '10+32=42' },\n    @{ Old = '0+14=14'; New = '7+38=45' },\n    @{ Old = '44-3=41'; New = '91-76=15' },\n    @{ Old = '50-13=37'; New = '95-10=85' },\n    @{ Old = '37+18=55'; New = '32+61=93' },\n    @{ Old = '4+53=57'; New = '23-12=11' },\n    @{ Old = '6+83=89'; New = '57+21=78' },\n    @{ Old = '4+6=10'; New = '13+83=96' },\n    @{ Old = '3+7=10'; New = '95-20=75' },\n    @{ Old = '37+47=84'; New = '36-32=4' },\n    @{ Old = '67-21=46'; New = '63-26=37' },\n    @{ Old = '97-37=60'; New = '98-47=51' },\n    @{ Old = '18+20=38'; New = '92-48=44' },\n    @{ Old = '89-6=83'; New = '98-31=67' },\n    @{ Old = '29+3=32'; New = '10+85=95' },\n    @{ Old = '38-11=27'; New = '87-38=49' },\n    @{ Old = '2+43=45'; New = '26+21=47' },\n    @{ Old = '16+72=88'; New = '70-30=40' },\n    @{ Old = '2+91=93'; New = '23+0=23' },\n    @{ Old = '52-28=24'; New = '26-12=14' },\n    @{ Old = '27+69=96'; New = '90-56=34' },\n    @{ Old = '45-38=7'; New = '40-4=36' },\n    @{ Old = '6+60=66'; New = '38-3=35' },\n    @{ Old = '58-18=40'; New = '50-36=14' },\n    @{ Old = '90-65=25'; New = '25+40=65' },\n    @{ Old = '64+28=92'; New = '51+30=81' },\n    @{ Old = '60+27=87'; New = '76-25=51' },\n    @{ Old = '80-53=27'; New = '40+29=69' },\n    @{ Old = '69-3=66'; New = '10+89=99' },\n    @{ Old = '1+39=40'; New = '79-24=55' },\n    @{ Old = '44+53=97'; New = '85+13=98' },\n    @{ Old = '36-34=2'; New = '98-61=37' },\n    @{ Old = '77+17=94'; New = '55+29=84' },\n    @{ Old = '11+22=33'; New = '1+53=54' },\n    @{ Old = '4+68=72'; New = '1+17=18' },\n    @{ Old = '51-34=17'; New = '1+30=31' },\n    @{ Old = '55-5=50'; New = '85-84=1' },\n    @{ Old = '63+30=93'; New = '74+14=88' },\n    @{ Old = '24+75=99'; New = '42+51=93' },\n    @{ Old = '44+52=96'; New = '11+73=84' },\n    @{ Old = '51-36=15'; New = '49+44=93' },\n    @{ Old = '11+76=87'; New = '71-53=18' },\n    @{ Old = '38+32=70'; New = '26+8=34' },\n    @{ Old = '68-47=21'; New = '76-32=44' },\n    @{ Old = '89-57=32'; New = '44+15=59' },\n    @{ Old = '48+36=84'; New = '62+31=93' },\n    @{ Old = '64-34=30'; New = '43+35=78' },\n    @{ Old = '73-67=6'; New = '70-7=63' },\n    @{ Old = '4+75=79'; New = '76+1=77' },\n    @{ Old = '57+42=99'; New = '13+33=46' },\n    @{ Old = '19+55=74'; New = '7+85=92' },\n    @{ Old = '28+18=46'; New = '4+11=15' },\n    @{ Old = '68-53=15'; New = '35+43=78' },\n    @{ Old = '96-27=69'; New = '17-2=15' },\n    @{ Old = '51-26=25'; New = '2+45=47' },\n    @{ Old = '19+54=73'; New = '37-17=20' },\n    @{ Old = '77-0=77'; New = '64-10=54' },\n    @{ Old = '80-72=8'; New = '41+45=86' },\n    @{ Old = '19-17=2'; New = '54-2=52' },\n    @{ Old = '23+76=99'; New = '68-57=11' },\n    @{ Old = '3+14=17'; New = '94-70=24' },\n    @{ Old = '41+22=63'; New = '33+48=81' },\n    @{ Old = '29+46=75'; New = '62-19=43' },\n    @{ Old = '43-28=15'; New = '71+27=98' },\n    @{ Old = '53+7=60'; New = '98-96=2' }\n)\n\n# wdReplace: 1 = wdReplaceOne (only the first match), 2 = wdReplaceAll\n# wdFindWrap: 1 = wdFindContinue (search from start of range, don't wrap past end)\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 1)\n    if (-not $found) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
